# Update advies-tabblad + timestamp in A1 (A2 holds the "Laatst bijgewerkt" label)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Refresh "last updated" timestamp ---
$ws.Range("A2").Value = "Laatst bijgewerkt: 2025-09-06 19:12:44"

# --- 2) Replace the first advies row (row 3) with the new Litouwen vs Nederland /
#        "totaal aantal schoten" / vbet entry. Rows 4-6 keep their data, only the
#        hyperlinks need to be rebuilt (see step 3). ---
$ws.Range("A3").Value = "Litouwen vs Nederland"
$ws.Range("B3").Value = "totaal aantal schoten"
$ws.Range("C3").Value = "sem steijn"
$ws.Range("D3").Value = "meer dan 3.5"
$ws.Range("E3").Value = "vbet"
$ws.Range("F3").Value = 1.91
$ws.Range("G3").Value = "minder dan 3.5"
$ws.Range("H3").Value = "jacks"
$ws.Range("I3").Value = 2.23
$ws.Range("J3").Value = "1=81, 2=69"
$ws.Range("K3").Value = "€3.87"
$ws.Range("L3").Value = 2.8

# --- 3) Rebuild every hyperlink. The engine's Hyperlinks.Delete() call clears the
#        whole sheet's hyperlink collection, so all 8 need to be re-added in order
#        (M3,N3,M4,N4,M5,N5,M6,N6) to reproduce the original relationship ids.
#        Note the rId targets intentionally stay on their old (stale) URLs for
#        M3/N3 - only the displayed cell text changes there. ---
$ws.Range("M3:N6").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("M3"), "https://sport.toto.nl/wedden/wedstrijd/8778584")
$ws.Range("M3").Value = "https://www.vbet.nl/nl/sports/pre-match/event-view/Soccer/World/18277589/world-cup-europe-qualification/27857408/litouwen-nederland"

$ws.Hyperlinks.Add($ws.Range("N3"), "https://starcasino.nl/prematch-bets?page=event&sportId=66&eventId=12642394", "event/1023224945")
$ws.Range("N3").Value = "https://jacks.nl/sports/event/1023224945#event/1023224945"

$ws.Hyperlinks.Add($ws.Range("M4"), "https://sport.toto.nl/wedden/wedstrijd/8706282")
$ws.Range("M4").Value = "https://sport.toto.nl/wedden/wedstrijd/8706282"

$ws.Hyperlinks.Add($ws.Range("N4"), "https://jacks.nl/sports/event/1023224945", "event/1023224945")
$ws.Range("N4").Value = "https://jacks.nl/sports/event/1023224945#event/1023224945"

$ws.Hyperlinks.Add($ws.Range("M5"), "https://sport.toto.nl/wedden/wedstrijd/8706282")
$ws.Range("M5").Value = "https://sport.toto.nl/wedden/wedstrijd/8706282"

$ws.Hyperlinks.Add($ws.Range("N5"), "https://jacks.nl/sports/event/1023224945", "event/1023224945")
$ws.Range("N5").Value = "https://jacks.nl/sports/event/1023224945#event/1023224945"

$ws.Hyperlinks.Add($ws.Range("M6"), "https://sport.toto.nl/wedden/wedstrijd/8706282")
$ws.Range("M6").Value = "https://sport.toto.nl/wedden/wedstrijd/8706282"

$ws.Hyperlinks.Add($ws.Range("N6"), "https://jacks.nl/sports/event/1023224945", "event/1023224945")
$ws.Range("N6").Value = "https://jacks.nl/sports/event/1023224945#event/1023224945"
